$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price (column D) and 1h volume/change (column E) values
# as scraped by the GitHub Actions job on Fri Sep  8 13:31:58 UTC 2023.
# For price values that look like plain decimal numbers, force the cell
# to Text format before assigning so Excel keeps the exact original
# string (including trailing zeros) instead of coercing it to a number,
# then restore the default "Normal" style so no formatting is changed.
$ws.Range("D2").Value = '25.887.59'
$ws.Range("E2").Value = '  +0.68%  '
$ws.Range("D3").Value = '1.630.66'
$ws.Range("E3").Value = '  +0.09%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.69'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.19%  '
$ws.Range("E6").Value = '  +0.35%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("E9").Value = '  +0.01%  '
$ws.Range("E10").Value = '  +0.81%  '
$ws.Range("E11").Value = '  -0.66%  '
$ws.Range("D12").Value = '1.856.02'
$ws.Range("E12").Value = '  +0.09%  '
$ws.Range("E13").Value = '  -0.49%  '
$ws.Range("D14").Value = '1.587.61'
$ws.Range("E14").Value = '  -1.83%  '
$ws.Range("E15").Value = '  -1.89%  '
$ws.Range("E16").Value = '  -0.37%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.79'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.14%  '
$ws.Range("D18").Value = '25.881.38'
$ws.Range("E18").Value = '  +0.63%  '
$ws.Range("E19").Value = '  -0.06%  '
$ws.Range("E20").Value = '  -1.31%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '192.78'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.58%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.95'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.52%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.25'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.31%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.79'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.58%  '
$ws.Range("E25").Value = '  -0.12%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '142.49'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.05%  '
$ws.Range("E27").Value = '  +1.92%  '
$ws.Range("E28").Value = '  +0.15%  '
$ws.Range("E29").Value = '  -0.10%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0499'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.11%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.31'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.52%  '
$ws.Range("E33").Value = '  -0.24%  '
$ws.Range("E34").Value = '  +0.05%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.41'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.59%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.900'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.30%  '
$ws.Range("D37").Value = '1.135.47'
$ws.Range("E37").Value = '  +0.00%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.550'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.52%  '
$ws.Range("E39").Value = '  -2.26%  '
$ws.Range("E40").Value = '  +0.53%  '
$ws.Range("E41").Value = '  -0.07%  '
$ws.Range("E42").Value = '  -1.00%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '99.22'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.28%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.801'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.24%  '
$ws.Range("D45").Value = '1.765.97'
$ws.Range("E45").Value = '  +0.14%  '
$ws.Range("E46").Value = '  +0.20%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '56.03'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.51%  '
$ws.Range("E48").Value = '  +4.77%  '
$ws.Range("E49").Value = '  +1.52%  '
$ws.Range("E50").Value = '  -0.91%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.62'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.05%  '
